# Applies the "Doing Updates for Financials" data edits to the TEAM financials sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEAM")

# ---- Income Statement: Net Income (FY2016 column, column F) 4400 -> 4100 ----
$ws.Range("F27").Value = 4100   # Net Income From Continuing Ops
$ws.Range("F33").Value = 4100   # Net Income
$ws.Range("F35").Value = 4100   # Net Income Applicable To Common Shares

# ---- Balance Sheet updates (column D = most recent period) ----
$ws.Range("D41").Value = 2820700   # Cash And Cash Equivalents
$ws.Range("D43").Value = 110500    # Net Receivables
$ws.Range("D45").Value = 52400     # Other Current Assets
$ws.Range("D46").Value = 1822000   # Total Current Assets

# Long Term Investments row: D gets a real number, E:J become "NA"
$ws.Range("D47").Value = 5200
$ws.Range("E47").Value = "NA"
$ws.Range("F47").Value = "NA"
$ws.Range("G47").Value = "NA"
$ws.Range("H47").Value = "NA"
$ws.Range("I47").Value = "NA"
$ws.Range("J47").Value = "NA"

$ws.Range("D48").Value = 103300    # Property Plant and Equipment
$ws.Range("D49").Value = 439100    # Goodwill
$ws.Range("D52").Value = 279600    # Other Assets
$ws.Range("D54").Value = 2421800   # Total Assets
$ws.Range("D57").Value = 130200    # Accounts Payable
$ws.Range("D59").Value = 427800    # Other Current Liabilities
$ws.Range("D60").Value = 444900    # Total Current Liabilities
$ws.Range("D62").Value = 469300    # Other Liabilities
$ws.Range("D66").Value = 1514500   # Total Liabilities
$ws.Range("D72").Value = 429100    # Retained Earnings
$ws.Range("D76").Value = 907300    # Total Stockholder Equity

# ---- Cash Flow Statement ----
$ws.Range("F81").Value = 4100      # Net Income

# Capital Expenditures row
$ws.Range("E91").Value = -15100
$ws.Range("G91").Value = -31800
$ws.Range("H91").Value = -8100
